$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Palavras-chave: ... amostra; pesquisa" -> append a closing period
#    as its own run (same run formatting as the run it follows).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute("amostra; pesquisa", $false, $false, $false, $false, $false, $true, 1, $false, "amostra; pesquisa.", 2)
Write-Output ("Keyword period appended: " + $found1)

# ---------------------------------------------------------------------
# 2) "Fonte: adaptado de BVS atenção primária à saúde" -> split into
#    two runs: "Fonte:" (bold) + " adaptado de BVS atenção primária à
#    saúde" (unchanged formatting).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("Fonte: adaptado de BVS atenção primária à saúde")
Write-Output ("Source caption found: " + $found2)

if ($found2) {
    $fonteRange = $d.Range($rng2.Start, $rng2.Start + 6)
    Write-Output ("Bolding text: [" + $fonteRange.Text + "]")
    $fonteRange.Font.Bold = 1
    $fonteRange.Font.BoldBi = 1
}

# ---------------------------------------------------------------------
# 3) Page margins: top 993 -> 1276 twips, header distance 708 -> 993
#    twips (values expressed in points for the PageSetup object).
# ---------------------------------------------------------------------
$ps = $d.PageSetup
$ps.TopMargin = 63.8
$ps.HeaderDistance = 49.65
Write-Output ("TopMargin now: " + $ps.TopMargin)
Write-Output ("HeaderDistance now: " + $ps.HeaderDistance)
